$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Discussed" column added to support online-review workflow.
# Use the workbook's default/unstyled cell format (new cells shouldn't
# inherit the row's customFormat style) and stamp every data row
# (2-71) in column T with the review status.
$col = $ws.Range("T2:T71")
$col.Style = "Normal"
$col.Value = "Discussed"

# Reflect the reviewer's current working selection after adding the column.
$col.Select()
